$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change for each data row (row index, then values for B,D,E,F,G,H,I,J,K,N in that order)
$cols = @("B","D","E","F","G","H","I","J","K","N")

$data = @(
    @(2, 0.07492757642405934, 0.1418663856174618, 0.1535630155413727, 1.850564403315317, 1.202286409497134, 1.140910202406872, 0.5122450281881115, 0.1979505855692381, 1.190017266316573, 1.537932231489793),
    @(3, 0.06601511057540677, 0.1385942412145624, 0.1501442596801645, 1.84070584341633, 1.193364651663103, 1.142526749673834, 0.5186530680228412, 0.1936152618165892, 1.086145941070612, 1.556008322597279),
    @(4, 0.06053932161339048, 0.1366390091034475, 0.1481237964890028, 1.835762985020736, 1.188717721887599, 1.144113671600081, 0.5228324154145678, 0.1910707320192984, 1.022743832006114, 1.567668189366525),
    @(5, 0.05830720056020766, 0.1358558362146596, 0.1473202248180279, 1.83402752926952, 1.187032498314892, 1.144909692670424, 0.524597065389651, 0.190063279428891, 0.9970015652793336, 1.572560722525038),
    @(6, 0.05793652224858192, 0.135726613906499, 0.1471879873070634, 1.833756185254259, 1.186765241654996, 1.14505088759023, 0.5248937999373915, 0.189897770765576, 0.9927328027855253, 1.573381641346372),
    @(7, 0.06050922100401124, 0.1366283918237627, 0.1481128791341, 1.835738451751084, 1.1886941511343, 1.144123802518507, 0.5228559649897466, 0.1910570259294815, 1.022396279416398, 1.567733600764369),
    @(8, 0.07185539379632644, 0.1407269977284926, 0.1523679009174792, 1.846934468509545, 1.199037367379262, 1.141344155175787, 0.5144037030127908, 0.1964313715297692, 1.154124643155285, 1.544048391869367),
    @(9, 0.09407029813274903, 0.1491901070444044, 0.1613367250750528, 1.877722025256048, 1.225941841818184, 1.140615868687149, 0.4997718658617494, 0.2079049200261949, 1.415425500740469, 1.502054921066994),
    @(10, 0.1103622987864838, 0.1556659751223606, 0.1683086823633246, 1.905762386032819, 1.249787279709039, 1.142970946467244, 0.4902071553266225, 0.2169097632605741, 1.609247835103588, 1.473917863808195),
    @(11, 0.1177659958613759, 0.1586677750059948, 0.1715639084066112, 1.919704025329651, 1.261530241178264, 1.144672479089508, 0.4861134648706944, 0.221132483339801, 1.697830332694366, 1.461707445148116),
    @(12, 0.1205683195935876, 0.159812477083392, 0.1728086226358627, 1.925154472847183, 1.266106461774768, 1.14540761065453, 0.4846003220175472, 0.2227497696681553, 1.731433561555775, 1.45716847603474),
    @(13, 0.1199648498491541, 0.1595655906710647, 0.1725400160416015, 1.923973005637038, 1.26511512479118, 1.14524524579906, 0.4849245559556579, 0.2224006461716357, 1.724193886598641, 1.458142248575481),
    @(14, 0.1179965718353913, 0.1587617905880876, 0.171666070665438, 1.920149005835, 1.261904132249413, 1.144731138053714, 0.4859882349348137, 0.2212651726683958, 1.700593710311239, 1.461332320204075),
    @(15, 0.1167907700439343, 0.1582704782527173, 0.1711323207279634, 1.917828986847525, 1.259954179353116, 1.144428062104197, 0.4866445951648739, 0.2205720387983661, 1.686145591942193, 1.463297384838075),
    @(16, 0.1098782793651623, 0.155470920734075, 0.168097629080016, 1.904875175098354, 1.24903792467191, 1.142872443491314, 0.490479869724858, 0.2166363467579941, 1.603467024422798, 1.474727710745171),
    @(17, 0.1056355967532738, 0.1537677617726274, 0.1662573697356891, 1.897232542855406, 1.242570979885897, 1.142079654862556, 0.4928986290322186, 0.2142543377067625, 1.55285167884773, 1.481890852191375),
    @(18, 0.1031946130442947, 0.1527934139723897, 0.165206774384032, 1.892948297901285, 1.238935619526728, 1.141682970046332, 0.4943140552757352, 0.2128961600696186, 1.52377784050941, 1.48606635938255),
    @(19, 0.1023680229589274, 0.152464422270242, 0.1648524130069973, 1.89151687714147, 1.237719200578312, 1.141558839964375, 0.494797452761869, 0.2124383448732203, 1.513940604888433, 1.487489635852072),
    @(20, 0.1060873119680394, 0.1539485215067913, 0.1664524536284659, 1.898034560726046, 1.243250672495464, 1.142157909315927, 0.4926386409975549, 0.2145066755220455, 1.558235758108253, 1.481122583066405),
    @(21, 0.1185747395584826, 0.158997669842762, 0.1719224430743367, 1.921267560856094, 1.262843761185138, 1.144879678302232, 0.4856748007935394, 0.2215981934023574, 1.707524055954195, 1.460393015055651),
    @(22, 0.1267283506254131, 0.1623441025705006, 0.1755675233917984, 1.937448930552378, 1.276403611319353, 1.147187836730609, 0.481339498127797, 0.2263392205252757, 1.805436141385883, 1.447339679548671),
    @(23, 0.1223773805260038, 0.1605538087148233, 0.1736156579871491, 1.928721207902058, 1.269097203450428, 1.145907436651896, 0.4836335552004591, 0.2237990971076727, 1.753147263271501, 1.454261192054245),
    @(24, 0.1058830971395111, 0.1538667850108624, 0.1663642331697801, 1.897671627056539, 1.242943126298911, 1.142122346385577, 0.4927561041851956, 0.2143925584942821, 1.55580153612317, 1.48146973898181),
    @(25, 0.0880651455685495, 0.1468551914804834, 0.1588433634824327, 1.868443887488468, 1.217950261391337, 1.140306164766969, 0.5035220723682414, 0.2047004087707052, 1.344414731869222, 1.512938462393761)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $row[$i + 1]
        $ws.Range("$col$r").Value = $val
    }
}

Write-Host "Updated $($data.Length) rows"
